$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '23.695.09'
$ws.Range('E2').Value = '  +0.87%  '

# Row 3
$ws.Range('D3').Value = '1.655.06'
$ws.Range('E3').Value = '  +0.94%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.43%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3813'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.76%  '

# Row 8
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.53'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.34%  '

# Row 9
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3613'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.27%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08197'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.30%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.230'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.08%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.22%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.498'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.54%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.403'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.47%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001230'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.73%  '

# Row 17
$ws.Range('D17').Value = '1.652.27'
$ws.Range('E17').Value = '  +1.07%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.26%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07011'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.00%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.829'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.50%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '

# Row 22
$ws.Range('E22').Value = '  +0.25%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.78'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.17%  '

# Row 24
$ws.Range('D24').Value = '23.692.59'
$ws.Range('E24').Value = '  +0.91%  '

# Row 25
$ws.Range('E25').Value = '  +0.29%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.033'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.51%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.50%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.91%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.210'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.75%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.19%  '

# Row 31
$ws.Range('D31').Value = '1.833.70'
$ws.Range('E31').Value = '  +0.97%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.018'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.32%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.216'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.48%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.96'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.69%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.056'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.11%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02797'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.30%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2515'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.79%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08776'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.03%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.081'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.63%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07027'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.99%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.98%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6993'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.334'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.30%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.22%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6507'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.29%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '

# Row 47
$ws.Range('E47').Value = '  +1.44%  '

# Row 48
$ws.Range('E48').Value = '  -0.23%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07909'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.78%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.181'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.02%  '
